$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.052.68"
$ws.Range("E2").Value = "  -6.39%  "
$ws.Range("D3").Value = "2.187.42"
$ws.Range("E3").Value = "  -7.14%  "
$ws.Range("E4").Value = "  -0.34%  "
$ws.Range("D5").Value = "'241.03"
$ws.Range("E5").Value = "  +0.45%  "
$ws.Range("D6").Value = "'0.620"
$ws.Range("E6").Value = "  -7.17%  "
$ws.Range("D7").Value = "'69.99"
$ws.Range("E7").Value = "  -4.51%  "
$ws.Range("E8").Value = "  -0.08%  "
$ws.Range("D9").Value = "'0.538"
$ws.Range("E9").Value = "  -11.40%  "
$ws.Range("D10").Value = "'36.27"
$ws.Range("E10").Value = "  +6.50%  "
$ws.Range("D11").Value = "'57.67"
$ws.Range("E11").Value = "  -4.91%  "
$ws.Range("D12").Value = "'0.0933"
$ws.Range("E12").Value = "  -8.20%  "
$ws.Range("E13").Value = "  -4.08%  "
$ws.Range("D14").Value = "'6.55"
$ws.Range("E14").Value = "  -9.22%  "
$ws.Range("D15").Value = "2.502.61"
$ws.Range("E15").Value = "  -7.69%  "
$ws.Range("D16").Value = "'14.53"
$ws.Range("E16").Value = "  -10.15%  "
$ws.Range("D17").Value = "'0.831"
$ws.Range("E17").Value = "  -8.49%  "
$ws.Range("D18").Value = "2.174.80"
$ws.Range("E18").Value = "  -7.75%  "
$ws.Range("D19").Value = "40.894.56"
$ws.Range("E19").Value = "  -6.91%  "
$ws.Range("D20").Value = "0.0₃0935"
$ws.Range("E20").Value = "  -8.84%  "
$ws.Range("D21").Value = "'73.51"
$ws.Range("E21").Value = "  -5.28%  "
$ws.Range("D22").Value = "'6.01"
$ws.Range("E22").Value = "  -7.86%  "
$ws.Range("D23").Value = "'230.40"
$ws.Range("E23").Value = "  -8.89%  "
$ws.Range("D24").Value = "'2.01"
$ws.Range("E24").Value = "  +8.31%  "
$ws.Range("D25").Value = "'0.999"
$ws.Range("E25").Value = "  -0.04%  "
$ws.Range("E26").Value = "  -4.98%  "
$ws.Range("D27").Value = "'2.39"
$ws.Range("E27").Value = "  -4.32%  "
$ws.Range("E28").Value = "  -5.05%  "
$ws.Range("D29").Value = "'9.67"
$ws.Range("E29").Value = "  -7.55%  "
$ws.Range("D30").Value = "'167.80"
$ws.Range("E30").Value = "  -4.86%  "
$ws.Range("D31").Value = "'20.15"
$ws.Range("E31").Value = "  -9.45%  "
$ws.Range("D32").Value = "'0.117"
$ws.Range("E32").Value = "  -8.84%  "
$ws.Range("E33").Value = "  -8.00%  "
$ws.Range("D34").Value = "'0.0697"
$ws.Range("E34").Value = "  -6.65%  "
$ws.Range("D35").Value = "'5.07"
$ws.Range("E35").Value = "  -4.75%  "
$ws.Range("D36").Value = "'4.54"
$ws.Range("E36").Value = "  -10.01%  "
$ws.Range("D37").Value = "'3.84"
$ws.Range("E37").Value = "  +1.48%  "
$ws.Range("D38").Value = "'23.33"
$ws.Range("E38").Value = "  +17.03%  "
$ws.Range("E39").Value = "  -6.81%  "
$ws.Range("D40").Value = "'0.0266"
$ws.Range("E40").Value = "  -3.26%  "
$ws.Range("D41").Value = "'5.75"
$ws.Range("E41").Value = "  -13.30%  "
$ws.Range("D42").Value = "'64.89"
$ws.Range("E42").Value = "  +0.10%  "
$ws.Range("D43").Value = "'4.84"
$ws.Range("E43").Value = "  -12.06%  "
$ws.Range("E44").Value = "  -4.91%  "
$ws.Range("D45").Value = "'0.190"
$ws.Range("E45").Value = "  -5.55%  "
$ws.Range("E46").Value = "  +0.24%  "
$ws.Range("D47").Value = "'0.0976"
$ws.Range("E47").Value = "  -8.09%  "
$ws.Range("D48").Value = "'4.51"
$ws.Range("E48").Value = "  +5.19%  "
$ws.Range("D49").Value = "'10.01"
$ws.Range("E49").Value = "  +4.99%  "
$ws.Range("E50").Value = "  -6.48%  "
$ws.Range("E51").Value = "  -6.35%  "
